# Auto-generated Excel COM-interop script to apply cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$TextValue
    )
    $cell = $ws.Range($CellRef)
    # Force text number format so numeric-looking strings (e.g. "1.00", "7.71")
    # are stored as literal text rather than being parsed into numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $TextValue
    # Reset style back to Normal so no stray style index is left on the cell.
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "56.598.99"
Set-TextValue "E2" "  -2.33%  "

# Row 3
Set-TextValue "D3" "2.987.05"
Set-TextValue "E3" "  -4.54%  "

# Row 4
Set-TextValue "E4" "  +0.04%  "

# Row 5
Set-TextValue "D5" "494.32"
Set-TextValue "E5" "  -5.47%  "

# Row 6
Set-TextValue "D6" "134.51"
Set-TextValue "E6" "  -0.48%  "

# Row 7
Set-TextValue "E7" "  -0.08%  "

# Row 8
Set-TextValue "D8" "2.985.28"
Set-TextValue "E8" "  -4.45%  "

# Row 9
Set-TextValue "E9" "  -4.60%  "

# Row 10
Set-TextValue "D10" "7.21"
Set-TextValue "E10" "  -0.39%  "

# Row 11
Set-TextValue "E11" "  -5.69%  "

# Row 12
Set-TextValue "D12" "0.352"
Set-TextValue "E12" "  -8.63%  "

# Row 13
Set-TextValue "E13" "  +0.38%  "

# Row 14
Set-TextValue "D14" "3.499.15"
Set-TextValue "E14" "  -4.55%  "

# Row 15
Set-TextValue "D15" "24.84"
Set-TextValue "E15" "  -2.43%  "

# Row 16
Set-TextValue "D16" "56.521.32"
Set-TextValue "E16" "  -2.29%  "

# Row 17
Set-TextValue "D17" "2.989.44"
Set-TextValue "E17" "  -4.52%  "

# Row 18
Set-TextValue "E18" "  -5.24%  "

# Row 19
Set-TextValue "D19" "5.83"
Set-TextValue "E19" "  +0.17%  "

# Row 20
Set-TextValue "D20" "12.27"
Set-TextValue "E20" "  -5.91%  "

# Row 21
Set-TextValue "D21" "7.71"

# Row 22
Set-TextValue "D22" "321.58"
Set-TextValue "E22" "  -6.79%  "

# Row 23
Set-TextValue "E23" "  -0.06%  "

# Row 24
Set-TextValue "D24" "0.462"
Set-TextValue "E24" "  -8.40%  "

# Row 25
Set-TextValue "D25" "61.15"
Set-TextValue "E25" "  -11.21%  "

# Row 26
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  +0.47%  "

# Row 27
Set-TextValue "E27" "  -2.63%  "

# Row 28
Set-TextValue "E28" "  -7.77%  "

# Row 29
Set-TextValue "E29" "  +0.02%  "

# Row 30
Set-TextValue "D30" "6.56"
Set-TextValue "E30" "  -3.30%  "

# Row 31
Set-TextValue "D31" "6.79"
Set-TextValue "E31" "  -1.08%  "

# Row 32
Set-TextValue "E32" "  -4.86%  "

# Row 33
Set-TextValue "E33" "  -7.74%  "

# Row 34
Set-TextValue "D34" "19.85"
Set-TextValue "E34" "  -8.01%  "

# Row 35
Set-TextValue "D35" "151.32"
Set-TextValue "E35" "  -4.54%  "

# Row 36
Set-TextValue "D36" "4.47"
Set-TextValue "E36" "  -6.79%  "

# Row 37
Set-TextValue "D37" "5.63"
Set-TextValue "E37" "  -9.09%  "

# Row 38
Set-TextValue "E38" "  -6.93%  "

# Row 39
Set-TextValue "D39" "0.0665"
Set-TextValue "E39" "  -3.97%  "

# Row 40
Set-TextValue "D40" "23.36"
Set-TextValue "E40" "  -7.28%  "

# Row 41
Set-TextValue "D41" "3.020.56"
Set-TextValue "E41" "  -4.42%  "

# Row 42
Set-TextValue "D42" "37.23"
Set-TextValue "E42" "  -7.55%  "

# Row 43
Set-TextValue "E43" "  +0.02%  "

# Row 44
Set-TextValue "E44" "  -5.55%  "

# Row 45
Set-TextValue "E45" "  -7.90%  "

# Row 46
Set-TextValue "D46" "1.42"
Set-TextValue "E46" "  -2.52%  "

# Row 47
Set-TextValue "B47" "Maker"
Set-TextValue "C47" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D47" "2.178.89"
Set-TextValue "E47" "  -3.26%  "

# Row 48
Set-TextValue "B48" "Filecoin"
Set-TextValue "C48" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D48" "3.54"
Set-TextValue "E48" "  -9.52%  "

# Row 49
Set-TextValue "B49" "dogwifhat"
Set-TextValue "C49" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D49" "1.92"
Set-TextValue "E49" "  +4.71%  "

# Row 50
Set-TextValue "B50" "VeChain"
Set-TextValue "C50" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D50" "0.0236"
Set-TextValue "E50" "  +0.86%  "

# Row 51
Set-TextValue "B51" "InjectiveProtocol"
Set-TextValue "C51" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D51" "19.21"
Set-TextValue "E51" "  -5.17%  "
